# "load site for systems"
# - Rename the existing sheet to "Systems"
# - Add a new, empty "AssetSystemAssociation" sheet after it
# - Add a new "Siteid" header column (H) to the Systems sheet
# - Leave the active sheet as "Systems" with a revised selection

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "Systems"
$ws.Range("H1").Value = "Siteid"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "AssetSystemAssociation"

$ws.Activate()
$ws.Range("N17").Select() | Out-Null
